$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Documento" value in G3
$ws.Range("G3").Value = 30668697042

# Update the "NumeroCalle" value in M3
$ws.Range("M3").Value = 309

# Reset the scrolled view so column A is shown first (removes topLeftCell="D1")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

$ws.Range("M4").Select()
